# Update team-specific transition-probability matrix (Buffalo_A) cell values
# to reflect newly recomputed probabilities (updated underlying counts).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 0.1762295081967213
$ws.Cells.Item(2, 3).Value2 = 0.6065573770491803
$ws.Cells.Item(2, 10).Value2 = 0.01639344262295082
$ws.Cells.Item(2, 16).Value2 = 0.1147540983606557
$ws.Cells.Item(2, 19).Value2 = 0.0860655737704918
$ws.Cells.Item(3, 2).Value2 = 0.006369426751592357
$ws.Cells.Item(3, 3).Value2 = 0.05732484076433121
$ws.Cells.Item(3, 10).Value2 = 0.03184713375796178
$ws.Cells.Item(3, 16).Value2 = 0.7197452229299363
$ws.Cells.Item(3, 19).Value2 = 0.1847133757961783
$ws.Cells.Item(4, 10).Value2 = 0.08163265306122448
$ws.Cells.Item(4, 15).Value2 = 0.02040816326530612
$ws.Cells.Item(4, 16).Value2 = 0.5918367346938775
$ws.Cells.Item(4, 19).Value2 = 0.3061224489795918
$ws.Cells.Item(6, 2).Value2 = 0.05970149253731343
$ws.Cells.Item(6, 4).Value2 = 0.009950248756218905
$ws.Cells.Item(6, 6).Value2 = 0.05970149253731343
$ws.Cells.Item(6, 10).Value2 = 0.2238805970149254
$ws.Cells.Item(6, 15).Value2 = 0.02985074626865672
$ws.Cells.Item(6, 17).Value2 = 0.1641791044776119
$ws.Cells.Item(6, 18).Value2 = 0.03482587064676617
$ws.Cells.Item(6, 19).Value2 = 0.417910447761194
$ws.Cells.Item(7, 2).Value2 = 0.08021390374331551
$ws.Cells.Item(7, 4).Value2 = 0.0160427807486631
$ws.Cells.Item(7, 6).Value2 = 0.0481283422459893
$ws.Cells.Item(7, 10).Value2 = 0.1283422459893048
$ws.Cells.Item(7, 15).Value2 = 0.0053475935828877
$ws.Cells.Item(7, 17).Value2 = 0.160427807486631
$ws.Cells.Item(7, 18).Value2 = 0.06417112299465241
$ws.Cells.Item(7, 19).Value2 = 0.4973262032085561
$ws.Cells.Item(8, 2).Value2 = 0.07130434782608695
$ws.Cells.Item(8, 4).Value2 = 0.01913043478260869
$ws.Cells.Item(8, 6).Value2 = 0.05565217391304348
$ws.Cells.Item(8, 10).Value2 = 0.1321739130434783
$ws.Cells.Item(8, 15).Value2 = 0.02086956521739131
$ws.Cells.Item(8, 17).Value2 = 0.1443478260869565
$ws.Cells.Item(8, 18).Value2 = 0.08695652173913043
$ws.Cells.Item(8, 19).Value2 = 0.4695652173913044
$ws.Cells.Item(9, 2).Value2 = 0.1184834123222749
$ws.Cells.Item(9, 4).Value2 = 0.01421800947867299
$ws.Cells.Item(9, 6).Value2 = 0.04265402843601896
$ws.Cells.Item(9, 10).Value2 = 0.1090047393364929
$ws.Cells.Item(9, 15).Value2 = 0.004739336492890996
$ws.Cells.Item(9, 17).Value2 = 0.1753554502369668
$ws.Cells.Item(9, 18).Value2 = 0.09004739336492891
$ws.Cells.Item(9, 19).Value2 = 0.4454976303317535
$ws.Cells.Item(10, 2).Value2 = 0.0857843137254902
$ws.Cells.Item(10, 4).Value2 = 0.02532679738562092
$ws.Cells.Item(10, 6).Value2 = 0.06535947712418301
$ws.Cells.Item(10, 10).Value2 = 0.1127450980392157
$ws.Cells.Item(10, 15).Value2 = 0.01879084967320261
$ws.Cells.Item(10, 17).Value2 = 0.2295751633986928
$ws.Cells.Item(10, 18).Value2 = 0.0727124183006536
$ws.Cells.Item(10, 19).Value2 = 0.3897058823529412
$ws.Cells.Item(11, 7).Value2 = 0.119047619047619
$ws.Cells.Item(11, 10).Value2 = 0.07823129251700681
$ws.Cells.Item(11, 11).Value2 = 0.2244897959183673
$ws.Cells.Item(11, 12).Value2 = 0.5612244897959183
$ws.Cells.Item(11, 19).Value2 = 0.01700680272108844
$ws.Cells.Item(12, 7).Value2 = 0.7529411764705882
$ws.Cells.Item(12, 10).Value2 = 0.1882352941176471
$ws.Cells.Item(12, 11).Value2 = 0.01176470588235294
$ws.Cells.Item(12, 12).Value2 = 0.02352941176470588
$ws.Cells.Item(12, 19).Value2 = 0.02352941176470588
$ws.Cells.Item(13, 7).Value2 = 0.813953488372093
$ws.Cells.Item(13, 10).Value2 = 0.1627906976744186
$ws.Cells.Item(13, 19).Value2 = 0.02325581395348837
$ws.Cells.Item(15, 6).Value2 = 0.0091324200913242
$ws.Cells.Item(15, 8).Value2 = 0.228310502283105
$ws.Cells.Item(15, 9).Value2 = 0.0730593607305936
$ws.Cells.Item(15, 10).Value2 = 0.2968036529680365
$ws.Cells.Item(15, 11).Value2 = 0.0639269406392694
$ws.Cells.Item(15, 13).Value2 = 0.0091324200913242
$ws.Cells.Item(15, 15).Value2 = 0.0593607305936073
$ws.Cells.Item(15, 19).Value2 = 0.2602739726027397
$ws.Cells.Item(16, 6).Value2 = 0.01775147928994083
$ws.Cells.Item(16, 8).Value2 = 0.242603550295858
$ws.Cells.Item(16, 9).Value2 = 0.08875739644970414
$ws.Cells.Item(16, 10).Value2 = 0.3609467455621302
$ws.Cells.Item(16, 11).Value2 = 0.09467455621301775
$ws.Cells.Item(16, 13).Value2 = 0.02366863905325444
$ws.Cells.Item(16, 15).Value2 = 0.03550295857988166
$ws.Cells.Item(16, 19).Value2 = 0.136094674556213
$ws.Cells.Item(17, 6).Value2 = 0.01746724890829694
$ws.Cells.Item(17, 8).Value2 = 0.2510917030567685
$ws.Cells.Item(17, 9).Value2 = 0.07641921397379912
$ws.Cells.Item(17, 10).Value2 = 0.388646288209607
$ws.Cells.Item(17, 11).Value2 = 0.08951965065502183
$ws.Cells.Item(17, 13).Value2 = 0.01310043668122271
$ws.Cells.Item(17, 15).Value2 = 0.06331877729257641
$ws.Cells.Item(17, 19).Value2 = 0.1004366812227074
$ws.Cells.Item(18, 6).Value2 = 0.02808988764044944
$ws.Cells.Item(18, 8).Value2 = 0.2134831460674157
$ws.Cells.Item(18, 9).Value2 = 0.1123595505617977
$ws.Cells.Item(18, 10).Value2 = 0.4044943820224719
$ws.Cells.Item(18, 11).Value2 = 0.1067415730337079
$ws.Cells.Item(18, 13).Value2 = 0.02247191011235955
$ws.Cells.Item(18, 15).Value2 = 0.03370786516853932
$ws.Cells.Item(18, 19).Value2 = 0.07865168539325842
$ws.Cells.Item(19, 6).Value2 = 0.01075268817204301
$ws.Cells.Item(19, 8).Value2 = 0.2365591397849462
$ws.Cells.Item(19, 9).Value2 = 0.08960573476702509
$ws.Cells.Item(19, 10).Value2 = 0.3505376344086021
$ws.Cells.Item(19, 11).Value2 = 0.0974910394265233
$ws.Cells.Item(19, 13).Value2 = 0.01935483870967742
$ws.Cells.Item(19, 14).Value2 = 0.0007168458781362007
$ws.Cells.Item(19, 15).Value2 = 0.06594982078853047
$ws.Cells.Item(19, 19).Value2 = 0.1290322580645161
